$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps the same text value ("HK_R_acc_G"); shared-string table bookkeeping
# is an internal serialization detail handled by the engine itself.
$ws.Range("A1").Value = "HK_R_acc_G"

$values = @(
    77.376599634369285,
    76.736745886654475,
    74.908592321755023,
    75.731261425959786,
    75.731261425959786,
    77.056672760511873,
    74.908592321755023,
    75.822669104204749,
    74.725776965265084,
    74.725776965265084,
    75.776965265082268,
    76.782449725776971,
    75.457038391224856,
    76.051188299817184,
    75.822669104204749,
    75.457038391224856,
    76.691042047531994,
    76.828153564899452,
    78.199268738574034,
    76.005484460694689,
    78.382084095063988,
    75.457038391224856,
    75.548446069469833,
    75.731261425959786,
    75.914076782449726,
    75.868372943327245,
    75.594149908592328,
    76.691042047531994,
    75.822669104204749,
    76.005484460694689,
    74.771480804387565,
    74.862888482632542,
    74.908592321755023,
    75.685557586837291,
    75.914076782449726,
    81.215722120658143,
    74.497257769652649,
    74.954296160877504,
    75.822669104204749,
    76.09689213893968,
    76.553930530164536,
    76.234003656307124,
    76.279707495429619,
    76.645338208409513,
    75.274223034734916,
    75.18281535648994,
    77.559414990859239,
    75.639853747714795
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
